$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 10 through 19 (they are being removed entirely)
$ws.Range("A10:F19").EntireRow.Delete()

# Update Student ID values in column A for rows 2-9
# (leading apostrophe forces these to stay text, matching the original cell type)
$ws.Range("A2").Value = "'211147"
$ws.Range("A3").Value = "'211086"
$ws.Range("A4").Value = "'201197"
$ws.Range("A5").Value = "'191258"
$ws.Range("A6").Value = "'200904"
$ws.Range("A7").Value = "'200003"
$ws.Range("A8").Value = "'191186"
$ws.Range("A9").Value = "'201990"
